# Insert a new data row before row 103 (shifts existing rows 103:174 down to 104:175)
# and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("103:103").Insert()

$ws.Range("A103").Value = 10
$ws.Range("B103").Value = "Vega Modelo de Temuco"
$ws.Range("C103").Value = "La Araucanía"
$ws.Range("D103").Value = 44957
$ws.Range("E103").Value = 9
$ws.Range("F103").Value = 100112031
$ws.Range("G103").Value = "Poroto verde"
$ws.Range("H103").Value = "Brío"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 100
$ws.Range("K103").Value = 1200
$ws.Range("L103").Value = 1200
$ws.Range("M103").Value = 1200
$ws.Range("N103").Value = "`$/kilo"
$ws.Range("O103").Value = "Región de La Araucanía"
$ws.Range("P103").Value = 1200
$ws.Range("Q103").Value = 1
$ws.Range("R103").Value = "Hortaliza"
